# Update gh-pages to output generated at 456a3b4
# "展览" sheet (sheet1) and "全部类型" sheet (sheet4) both carry the same
# event rows; the numeric stat columns (registration/interest counts) for
# the two events were refreshed.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 5310
$wsExhibit.Range("F8").Value = 110

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value = 5310
$wsAll.Range("F11").Value = 110
